$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

$ws.Cells.Item(2, 1).Value = "이베스트스팩6호"
$ws.Cells.Item(2, 2).Value = "2024.06.27~06.28"
$ws.Cells.Item(2, 3).Value = "2,000~2,000"
$ws.Cells.Item(2, 4).Value = "-"
$ws.Cells.Item(2, 5).Value = 8000
$ws.Cells.Item(2, 6).Value = "이베스트 투자증권"

$ws.Cells.Item(3, 1).Value = "이엔셀"
$ws.Cells.Item(3, 2).Value = "2024.06.17~06.21"
$ws.Cells.Item(3, 3).Value = "13,600~15,300"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = 21308
$ws.Cells.Item(3, 6).Value = "NH투자증권"

$ws.Cells.Item(4, 1).Value = "하스"
$ws.Cells.Item(4, 2).Value = "2024.06.13~06.19"
$ws.Cells.Item(4, 3).Value = "9,000~12,000"
$ws.Cells.Item(4, 4).Value = "-"
$ws.Cells.Item(4, 5).Value = 16290
$ws.Cells.Item(4, 6).Value = "삼성증권"

$ws.Cells.Item(5, 1).Value = "에이치브이엠(구.한국진공야금)"
$ws.Cells.Item(5, 2).Value = "2024.06.11~06.17"
$ws.Cells.Item(5, 3).Value = "11,000~14,200"
$ws.Cells.Item(5, 4).Value = "-"
$ws.Cells.Item(5, 5).Value = 26400
$ws.Cells.Item(5, 6).Value = "NH투자증권"

$ws.Cells.Item(6, 1).Value = "이노스페이스"
$ws.Cells.Item(6, 2).Value = "2024.06.11~06.17"
$ws.Cells.Item(6, 3).Value = "36,400~43,300"
$ws.Cells.Item(6, 4).Value = "-"
$ws.Cells.Item(6, 5).Value = 48412
$ws.Cells.Item(6, 6).Value = "미래에셋증권,신한투자증권"

$ws.Cells.Item(7, 1).Value = "한국스팩15호"
$ws.Cells.Item(7, 2).Value = "2024.06.10~06.11"
$ws.Cells.Item(7, 3).Value = "2,000~2,000"
$ws.Cells.Item(7, 4).Value = "-"
$ws.Cells.Item(7, 5).Value = 12500
$ws.Cells.Item(7, 6).Value = "한국투자증권"

$ws.Cells.Item(8, 1).Value = "하이젠알앤엠"
$ws.Cells.Item(8, 2).Value = "2024.06.07~06.13"
$ws.Cells.Item(8, 3).Value = "4,500~5,500"
$ws.Cells.Item(8, 4).Value = "-"
$ws.Cells.Item(8, 5).Value = 15300
$ws.Cells.Item(8, 6).Value = "한국투자증권"

$ws.Cells.Item(9, 1).Value = "미래에셋비전스팩6호"
$ws.Cells.Item(9, 2).Value = "2024.06.05~06.07"
$ws.Cells.Item(9, 3).Value = "2,000~2,000"
$ws.Cells.Item(9, 4).Value = "-"
$ws.Cells.Item(9, 5).Value = 12900
$ws.Cells.Item(9, 6).Value = "미래에셋증권"

$ws.Cells.Item(10, 1).Value = "KB스팩29호"
$ws.Cells.Item(10, 2).Value = "2024.06.04~06.05"
$ws.Cells.Item(10, 3).Value = "2,000~2,000"
$ws.Cells.Item(10, 4).Value = "-"
$ws.Cells.Item(10, 5).Value = 12000
$ws.Cells.Item(10, 6).Value = "KB증권"

$ws.Cells.Item(11, 1).Value = "에이치엠씨아이비스팩7호"
$ws.Cells.Item(11, 2).Value = "2024.06.04~06.05"
$ws.Cells.Item(11, 3).Value = "2,000~2,000"
$ws.Cells.Item(11, 4).Value = "-"
$ws.Cells.Item(11, 5).Value = 14000
$ws.Cells.Item(11, 6).Value = "현대차증권"

$ws.Cells.Item(12, 1).Value = "에스오에스랩"
$ws.Cells.Item(12, 2).Value = "2024.06.03~06.10"
$ws.Cells.Item(12, 3).Value = "7,500~9,000"
$ws.Cells.Item(12, 4).Value = "-"
$ws.Cells.Item(12, 5).Value = 15000
$ws.Cells.Item(12, 6).Value = "한국투자증권"

$ws.Cells.Item(13, 1).Value = "미래에셋비전스팩5호"
$ws.Cells.Item(13, 2).Value = "2024.06.03~06.04"
$ws.Cells.Item(13, 3).Value = "2,000~2,000"
$ws.Cells.Item(13, 4).Value = "-"
$ws.Cells.Item(13, 5).Value = 9500
$ws.Cells.Item(13, 6).Value = "미래에셋증권"

$ws.Cells.Item(14, 1).Value = "한국스팩14호"
$ws.Cells.Item(14, 2).Value = "2024.06.03~06.04"
$ws.Cells.Item(14, 3).Value = "2,000~2,000"
$ws.Cells.Item(14, 4).Value = "-"
$ws.Cells.Item(14, 5).Value = 8000
$ws.Cells.Item(14, 6).Value = "한국투자증권"

$ws.Cells.Item(15, 1).Value = "엑셀세라퓨틱스"
$ws.Cells.Item(15, 2).Value = "2024.06.03~06.10"
$ws.Cells.Item(15, 3).Value = "6,200~7,700"
$ws.Cells.Item(15, 4).Value = "-"
$ws.Cells.Item(15, 5).Value = 10032
$ws.Cells.Item(15, 6).Value = "대신증권"

$ws.Cells.Item(16, 1).Value = "시프트업"
$ws.Cells.Item(16, 2).Value = "2024.06.03~06.13"
$ws.Cells.Item(16, 3).Value = "47,000~60,000"
$ws.Cells.Item(16, 4).Value = "-"
$ws.Cells.Item(16, 5).Value = 340750
$ws.Cells.Item(16, 6).Value = "한국투자증권,NH투자증권,신한투자증권"

$ws.Cells.Item(17, 1).Value = "이노그리드"
$ws.Cells.Item(17, 2).Value = "2024.05.31~06.07"
$ws.Cells.Item(17, 3).Value = "29,000~35,000"
$ws.Cells.Item(17, 4).Value = "-"
$ws.Cells.Item(17, 5).Value = 17400
$ws.Cells.Item(17, 6).Value = "한국투자증권"

$ws.Cells.Item(18, 1).Value = "한중엔시에스"
$ws.Cells.Item(18, 2).Value = "2024.05.30~06.05"
$ws.Cells.Item(18, 3).Value = "20,000~23,500"
$ws.Cells.Item(18, 4).Value = "-"
$ws.Cells.Item(18, 5).Value = 32000
$ws.Cells.Item(18, 6).Value = "IBK투자증권"

$ws.Cells.Item(19, 1).Value = "디비금융스팩12호"
$ws.Cells.Item(19, 2).Value = "2024.05.28~05.29"
$ws.Cells.Item(19, 3).Value = "2,000~2,000"
$ws.Cells.Item(19, 4).Value = "-"
$ws.Cells.Item(19, 5).Value = 10000
$ws.Cells.Item(19, 6).Value = "DB금융투자"

$ws.Cells.Item(20, 1).Value = "씨어스테크놀로지"
$ws.Cells.Item(20, 2).Value = "2024.05.27~05.31"
$ws.Cells.Item(20, 3).Value = "10,500~14,000"
$ws.Cells.Item(20, 4).Value = "-"
$ws.Cells.Item(20, 5).Value = 13650
$ws.Cells.Item(20, 6).Value = "한국투자증권"

$ws.Cells.Item(21, 1).Value = "라메디텍"
$ws.Cells.Item(21, 2).Value = "2024.05.27~05.31"
$ws.Cells.Item(21, 3).Value = "10,400~12,700"
$ws.Cells.Item(21, 4).Value = "-"
$ws.Cells.Item(21, 5).Value = 13499
$ws.Cells.Item(21, 6).Value = "대신증권"

Write-Output "done"